$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$newPara = $titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaLabel = "Meta description"
$metaRest = ": Get ready to play Baboon to the Moon, a unique and exciting slot machine with a chance to win up to 5000 times your bet. Try it for free now!"

$metaXml = "<w:p xmlns:w='$wNs'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>$metaLabel</w:t></w:r><w:r><w:t>$metaRest</w:t></w:r></w:p>"
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicate bold title paragraph
#    ("Play Baboon to the Moon Free - Exciting Slot Machine") and replace the
#    text of the remaining italic paragraph with the new DALL-E image prompt.
# ---------------------------------------------------------------------------
$dupTitleText = "Play Baboon to the Moon Free - Exciting Slot Machine"

$dupPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq $dupTitleText) {
        $dupPara = $para
        break
    }
}
$dupPara.Range.Delete()

$promptPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptText = "Prompt: DALLE, we would like you to create a feature image that captures the fun and adventurous theme of Baboon to the Moon. The image should be in a cartoon style and should prominently feature a happy Maya warrior wearing glasses as the main character. The image should also include symbols or elements related to the game, such as the monkey, moons, and playing cards. The overall vibe should be exciting and playful, enticing players to give this unique slot game a try."

$promptXml = "<w:p xmlns:w='$wNs'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$promptText</w:t></w:r></w:p>"
$promptPara.Range.InsertXML($promptXml)

Write-Host "Done."
